# DOMA-2542 Localization for Excel template (ticket_report_status_executor)
#
# The template uses a "{d.tickets[i + 1].<field>}" placeholder style for the
# second (overflow) ticket row. Remove the spaces around the "+ 1" so the
# templating engine renders it as "{d.tickets[i+1].<field>}", matching the
# already-tight "{d.tickets[i].<field>}" style used one row above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fields = @(
    "categoryClassifier",
    "address",
    "processing",
    "completed",
    "canceled",
    "deferred",
    "closed",
    "new_or_reopened"
)

for ($col = 1; $col -le $fields.Length; $col++) {
    $cell = $ws.Cells.Item(3, $col)
    $field = $fields[$col - 1]
    $cell.Value = "{d.tickets[i+1]." + $field + "}"
}
